$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation is inserted as the new row 2, pushing every
# existing data row (2..11) down by one (so old row 2 becomes row 3, ...,
# old row 11 becomes row 12). We do this by copying whole rows bottom-up
# (copy carries values + number formats/styles, e.g. the date style on
# column D) so nothing needs to be re-created from scratch.
for ($r = 11; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# Populate the newly freed row 2 with the new record's values (everything
# besides D/J/K/L/M/P stays identical to the row below it, which the copy
# above already took care of).
$ws.Range("D2").Value = 44648
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 6500
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 6750
$ws.Range("P2").Value = 112
